$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph: "Data analysis protocol" (was split across several
#    runs with de-CH language + spell-check markers). Rebuild it as a
#    single en-US run with no proofErr markers.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$nextPara = $d.Paragraphs.Item(2)

# Create a brand-new (clean) paragraph right before paragraph 2, inheriting
# paragraph 2's (en-US) paragraph-mark formatting.
$nextPara.Range.InsertParagraphBefore()

# The old title paragraph (with de-CH / proofErr junk) is still paragraph 1;
# the fresh, empty paragraph is now paragraph 2. Delete the old one.
$d.Paragraphs.Item(1).Range.Delete()

# Paragraph 1 is now the clean, empty, en-US paragraph. Fill it in.
$d.Paragraphs.Item(1).Range.InsertAfter("Data analysis protocol")

# ---------------------------------------------------------------------
# 2) Append four paragraphs after the final paragraph of the document:
#      (blank)
#      Land cover data:
#      (blank)
#      (blank)
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphAfter()

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$p.Range.InsertParagraphAfter()

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$p.Range.InsertBefore("Land cover data:")

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$p.Range.InsertParagraphAfter()

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$p.Range.InsertParagraphAfter()
